$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Cohort" column from the CasesTab (B2) Neo4j query and fix the
# trailing comma left behind on the preceding "Response to Treatment" line.
$newB2 = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
WHERE diag.primary_disease_site IN [''Lymph Node'']
RETURN  coalesce(c.case_id, '''') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '''') AS `Study Code` ,
        coalesce(s.clinical_study_type, '''') AS  `Study Type`,
        coalesce(demo.breed, '''') AS Breed ,
        coalesce(diag.disease_term, '''') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '''') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '''') AS Age ,
        coalesce(demo.sex, '''') AS Sex ,
        coalesce(demo.neutered_indicator, '''') AS `Neutered Status`,
        coalesce(demo.weight, '''') AS `Weight (kg)`,
        coalesce(diag.best_response, '''') AS `Response to Treatment`'

$ws.Range("B2").Value = $newB2

# Reflect the row shrinking by one wrapped line now that the Cohort clause is gone.
$ws.Rows(2).RowHeight = 244.8
$ws.Rows(3).RowHeight = 244.8
$ws.Rows(4).RowHeight = 244.8

# Leave the selection on the edited cell, like the author did before saving.
$ws.Range("B2").Select() | Out-Null
